$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Arkusz2")
$xvals = @(360,608,845,1056,1228,1393,1591,1860,2129)
$yvals = @(777,1974,3711,5089,7872,12466,17928,27638,40653)
for ($i = 0; $i -lt 9; $i++) {
    $ws2.Cells.Item(1, $i+1).Value = $xvals[$i]
    $ws2.Cells.Item(2, $i+1).Value = $yvals[$i]
}
$chartObj = $ws2.ChartObjects().Add(100, 100, 400, 300)
$chart = $chartObj.Chart
$chart.ChartType = 74
$chart.SeriesCollection.NewSeries()
$ser = $chart.SeriesCollection(1)
$ser.XValues = "=Arkusz2!`$A`$1:`$I`$1"
$ser.Values = "=Arkusz2!`$A`$2:`$I`$2"
$ser.Smooth = $true
$tl = $ser.Trendlines().Add(3)
$tl.Order = 3
$tl.DisplayRSquared = $true
$tl.DisplayEquation = $true
